# Fix a typo/text-correction in the "ZoneTexte 39" text box:
# Replace the run span describing experimental levels with the corrected wording,
# merging the previously split/misspelled runs into one corrected run.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$textShape = $s.Shapes.Item("ZoneTexte 39")
$tr = $textShape.TextFrame.TextRange
$sub = $tr.Characters(62, 362)
$sub.Text = " with the main workflow results, optionally gathering results obtained at different experimental (different molecular levels, different time points, different pre-exposure histories, …) extended with additional columns coding for the biological annotation of items and optionally for the experimental. Some lines of the workflow results can be replicated for items having more than one annotation "

# Reposition the "Image 47" picture slightly (offset change only, same size).
$pic = $s.Shapes.Item("Image 47")
$pic.Left = 512.838623046875
$pic.Top = 120.36150360107422
